$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.462.32"
$ws.Range("E2").Value = "  +1.81%  "
$ws.Range("D3").Value = "1.635.70"
$ws.Range("E3").Value = "  +3.17%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.54"
$ws.Range("E5").Value = "  +2.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3768"
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.85"
$ws.Range("E8").Value = "  +3.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3669"
$ws.Range("E9").Value = "  +2.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.274"
$ws.Range("E10").Value = "  +3.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08189"
$ws.Range("E11").Value = "  +2.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("E13").Value = "  +4.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.652"
$ws.Range("E14").Value = "  +2.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001281"
$ws.Range("E15").Value = "  +3.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.458"
$ws.Range("E16").Value = "  +1.70%  "
$ws.Range("D17").Value = "1.635.11"
$ws.Range("E17").Value = "  +2.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.69"
$ws.Range("E18").Value = "  +2.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06949"
$ws.Range("E19").Value = "  +3.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.36"
$ws.Range("E20").Value = "  +2.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.575"
$ws.Range("E21").Value = "  +2.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("D23").Value = "23.473.22"
$ws.Range("E23").Value = "  +1.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.91"
$ws.Range("E24").Value = "  +1.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.117"
$ws.Range("E25").Value = "  +9.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.407"
$ws.Range("E26").Value = "  +1.47%  "
$ws.Range("E27").Value = "  +3.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.30"
$ws.Range("E28").Value = "  +2.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.326"
$ws.Range("E29").Value = "  +2.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "136.14"
$ws.Range("E30").Value = "  +3.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.419"
$ws.Range("E31").Value = "  +3.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.803"
$ws.Range("E32").Value = "  +3.40%  "
$ws.Range("D33").Value = "1.816.98"
$ws.Range("E33").Value = "  +3.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9753"
$ws.Range("E34").Value = "  +3.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02814"
$ws.Range("E35").Value = "  +5.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.45"
$ws.Range("E36").Value = "  +4.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.07432"
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.227"
$ws.Range("E38").Value = "  +2.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2540"
$ws.Range("E39").Value = "  +2.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08834"
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.396"
$ws.Range("E41").Value = "  +3.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7157"
$ws.Range("E42").Value = "  +2.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.64"
$ws.Range("E43").Value = "  +4.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.15"
$ws.Range("E44").Value = "  +8.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6620"
$ws.Range("E45").Value = "  +3.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.358"
$ws.Range("E46").Value = "  +4.22%  "
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.044"
$ws.Range("E48").Value = "  +1.64%  "
$ws.Range("E49").Value = "  +2.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "130.89"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.216"
$ws.Range("E51").Value = "  +1.62%  "
